$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 108, shifting the existing data (rows 108-203)
# down to rows 109-204.
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new record.
$ws.Cells.Item(108, 1).Value2 = 3
$ws.Cells.Item(108, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(108, 3).Value2 = "Coquimbo"
$ws.Cells.Item(108, 4).Value2 = 44827
$ws.Cells.Item(108, 5).Value2 = 5
$ws.Cells.Item(108, 6).Value2 = 100112026
$ws.Cells.Item(108, 7).Value2 = "Haba"
$ws.Cells.Item(108, 8).Value2 = "Sin especificar"
$ws.Cells.Item(108, 9).Value2 = "Primera"
$ws.Cells.Item(108, 10).Value2 = 95
$ws.Cells.Item(108, 11).Value2 = 10000
$ws.Cells.Item(108, 12).Value2 = 11000
$ws.Cells.Item(108, 13).Value2 = 10474
$ws.Cells.Item(108, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(108, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(108, 16).Value2 = 419
$ws.Cells.Item(108, 17).Value2 = 25
$ws.Cells.Item(108, 18).Value2 = "Hortaliza"
